$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country order (several countries in the live feed moved up in rank)
# --- and refresh country statistics, per "Update countries & provincias Spain".


# Row 1
$ws.Range("A1").Value = 'Datos actualizados a 12 de Abril de 2020 a las 14:52'

# Row 15
$ws.Range("E15").Value = 12119
$ws.Range("G15").Value = 45
$ws.Range("H15").Value = 1081

# Row 17
$ws.Range("B17").Value = 20984
$ws.Range("C17").Value = 22
$ws.Range("E17").Value = 19670

# Row 23
$ws.Range("B23").Value = 10483
$ws.Range("C23").Value = 332
$ws.Range("E23").Value = 9203
$ws.Range("F23").Value = 839
$ws.Range("G23").Value = 12
$ws.Range("H23").Value = 899

# Row 39
$ws.Range("A39").Value = 'Arabia Saudita'
$ws.Range("B39").Value = 4462
$ws.Range("C39").Value = 429
$ws.Range("D39").Value = 761
$ws.Range("E39").Value = 3642
$ws.Range("F39").Value = 67
$ws.Range("G39").Value = 7
$ws.Range("H39").Value = 59

# Row 40
$ws.Range("A40").Value = 'Indonesia'
$ws.Range("B40").Value = 4241
$ws.Range("C40").Value = 399
$ws.Range("D40").Value = 359
$ws.Range("E40").Value = 3509
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 46
$ws.Range("H40").Value = 373

# Row 41
$ws.Range("A41").Value = 'Mexico'
$ws.Range("B41").Value = 4219
$ws.Range("C41").Value = 375
$ws.Range("D41").Value = 1772
$ws.Range("E41").Value = 2174
$ws.Range("F41").Value = 89
$ws.Range("G41").Value = 40
$ws.Range("H41").Value = 273

# Row 55
$ws.Range("E55").Value = 1714
$ws.Range("G55").Value = 5
$ws.Range("H55").Value = 98

# Row 69
$ws.Range("A69").Value = 'Barein'
$ws.Range("B69").Value = 1087
$ws.Range("C69").Value = 47
$ws.Range("D69").Value = 557
$ws.Range("E69").Value = 524
$ws.Range("F69").Value = 3
$ws.Range("H69").Value = 6

# Row 70
$ws.Range("A70").Value = 'Azerbaiyan'
$ws.Range("B70").Value = 1058
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 200
$ws.Range("E70").Value = 847
$ws.Range("F70").Value = 27
$ws.Range("H70").Value = 11

# Row 71
$ws.Range("A71").Value = 'Lituania'
$ws.Range("B71").Value = 1053
$ws.Range("C71").Value = 27
$ws.Range("D71").Value = 97
$ws.Range("E71").Value = 933
$ws.Range("F71").Value = 14
$ws.Range("H71").Value = 23

# Row 109
$ws.Range("D109").Value = 58
$ws.Range("E109").Value = 208

# Row 117
$ws.Range("A117").Value = 'Kenia'
$ws.Range("B117").Value = 197
$ws.Range("C117").Value = 6
$ws.Range("D117").Value = 25
$ws.Range("E117").Value = 164
$ws.Range("F117").Value = 2
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 8

# Row 118
$ws.Range("A118").Value = 'Mayotte'
$ws.Range("B118").Value = 196
$ws.Range("D118").Value = 59
$ws.Range("E118").Value = 134
$ws.Range("F118").Value = 3
$ws.Range("H118").Value = 3

# Row 132
$ws.Range("B132").Value = 106
$ws.Range("C132").Value = 4
$ws.Range("D132").Value = 20
$ws.Range("E132").Value = 86

# Row 133
$ws.Range("A133").Value = 'Mali'
$ws.Range("B133").Value = 105
$ws.Range("C133").Value = 18
$ws.Range("D133").Value = 22
$ws.Range("E133").Value = 74
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 2
$ws.Range("H133").Value = 9

# Row 134
$ws.Range("A134").Value = 'Monaco'
$ws.Range("D134").Value = 5
$ws.Range("E134").Value = 86
$ws.Range("F134").Value = 4
$ws.Range("H134").Value = 1

# Row 135
$ws.Range("A135").Value = 'Aruba'
$ws.Range("B135").Value = 92
$ws.Range("D135").Value = 29
$ws.Range("E135").Value = 63
$ws.Range("H135").Value = 0

# Row 146
$ws.Range("A146").Value = 'Liberia'
$ws.Range("C146").Value = 2
$ws.Range("D146").Value = 3
$ws.Range("E146").Value = 42
$ws.Range("F146").Value = 0
$ws.Range("H146").Value = 5

# Row 147
$ws.Range("A147").Value = 'San Martin (Parte Holandesa)'
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 5
$ws.Range("E147").Value = 36
$ws.Range("H147").Value = 9

# Row 148
$ws.Range("A148").Value = 'Bermudas'
$ws.Range("B148").Value = 50
$ws.Range("C148").Value = 2
$ws.Range("D148").Value = 28
$ws.Range("E148").Value = 18
$ws.Range("F148").Value = 2
$ws.Range("H148").Value = 4

# Row 149
$ws.Range("A149").Value = 'Gabon'
$ws.Range("B149").Value = 49
$ws.Range("C149").Value = 3
$ws.Range("D149").Value = 1
$ws.Range("E149").Value = 47
$ws.Range("H149").Value = 1

# Row 167
$ws.Range("A167").Value = 'Guinea Ecuatorial'
$ws.Range("C167").Value = 3
$ws.Range("D167").Value = 3
$ws.Range("F167").Value = 0
$ws.Range("H167").Value = 0

# Row 168
$ws.Range("A168").Value = 'Somalia'
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 2
$ws.Range("F168").Value = 2
$ws.Range("H168").Value = 1

# Row 187
$ws.Range("A187").Value = 'Nepal'
$ws.Range("C187").Value = 3

# Row 188
$ws.Range("A188").Value = 'San Vicente y las Granadinas'
$ws.Range("C188").Value = 0

# Row 210
$ws.Range("A210").Value = 'Anguila'
$ws.Range("C210").Value = 0

# Row 211
$ws.Range("A211").Value = 'Bonaire, San Eustaquio y Saba'
$ws.Range("C211").Value = 1

# Row 215
$ws.Range("A215").Value = 'San Pedro y Miquelon'

# Row 216
$ws.Range("A216").Value = 'Yemen'
